$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C12").Value = "showroom"
$ws.Range("D12").Value = "Roiis"
$ws.Range("E12").Value = "DONE"

$ws.Range("C13").Value = "dealer"
$ws.Range("D13").Value = "Roiis"
$ws.Range("E13").Value = "DONE"

$ws.Range("C14").Value = "listing"
$ws.Range("D14").Value = "Roiis"
$ws.Range("E14").Value = "DONE"

$ws.Range("E15").Select() | Out-Null
